{"js": "// Insert three new bullet paragraphs after the \"Data Engineering and\n// Infrastructure Architecture\" paragraph (under the PARTNER - Siege\n// Analytics heading), before the existing \"Architect enterprise-scale...\"\n// bullet.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text.trim() === \"Data Engineering and Infrastructure Architecture\"\n);\n\nif (!target) {\n  throw new Error(\"Could not find anchor paragraph 'Data Engineering and Infrastructure Architecture'\");\n}\n\nconst newBullets = [\n  \"\\u2022 Architected data infrastructure processing 15+ billion voter records to support meta-analytical voter file corrections\",\n  \"\\u2022 Built scalable ETL pipelines enabling analysis of 50,000+ electoral boundaries across all levels of government\",\n  \"\\u2022 Developed Python boundary estimation algorithm that reduced mapping costs by 75% for 200+ organizations\"\n];\n\n// Insert in order, each time right after the anchor paragraph so that the\n// bullets end up in the same order as in the diff.\nlet insertAfter = target;\nfor (const text of newBullets) {\n  insertAfter = insertAfter.insertParagraph(text, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Insert three new bullet paragraphs after the \"Data Engineering and\n# Infrastructure Architecture\" paragraph (under the PARTNER - Siege\n# Analytics heading), before the existing \"Architect enterprise-scale...\"\n# bullet.\n\n$d = $word.ActiveDocument\n\n$targetIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text.Trim() -eq \"Data Engineering and Infrastructure Architecture\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find anchor paragraph 'Data Engineering and Infrastructure Architecture'\"\n}\n\n$bulletChar = [char]0x2022\n\n$bulletLines = @(\n    \"$bulletChar Architected data infrastructure processing 15+ billion voter records to support meta-analytical voter file corrections\",\n    \"$bulletChar Built scalable ETL pipelines enabling analysis of 50,000+ electoral boundaries across all levels of government\",\n    \"$bulletChar Developed Python boundary estimation algorithm that reduced mapping costs by 75% for 200+ organizations\"\n)\n\n$insertIndex = $targetIndex\nforeach ($line in $bulletLines) {\n    $anchorPara = $d.Paragraphs.Item($insertIndex)\n    $anchorPara.Range.InsertParagraphAfter()\n    $insertIndex = $insertIndex + 1\n    $newPara = $d.Paragraphs.Item($insertIndex)\n    $newPara.Range.Text = $line\n}\n"}
